# "nuevos cambios desde oficina" - add a new data row (row 9) below the
# existing summary table, mirroring the pattern of the other detail rows
# (amount / sueldo-base / sub. al empleo efectivo columns A, C, E) plus a
# total formula in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of figures.
$ws.Range("A9").Value = 2038.91
$ws.Range("C9").Value = 7865.64
$ws.Range("E9").Value = 7865.64

# Row total, same shape as the existing row-7 total formula.
$ws.Range("G9").Formula = "=SUM(A9:F9)"

# Column A needed a touch more room for the new figures.
$ws.Columns.Item(1).ColumnWidth = 7.14

# Leave the new row selected, like it was just filled in.
$ws.Range("A9:C9").Select()
